$d = $word.ActiveDocument
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$anchorPara = $d.Paragraphs.Last
$insertRange = $anchorPara.Range
$insertRange.Collapse(1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">23/11/2023</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Transformación</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Guardar siempre una copia de los datos en raw.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Tener en cuenta los siguientes pasos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Normalización: qué datos son importantes y cómo se deben guardar</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Eliminar duplicados</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Verificación: comprobaciones automatizadas de que la información existe. Generar alarmas si alguno de los pasos falla</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Clasificaición: agrupar y clasificar los datos en bruto</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Un proceso de transformación</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Debe: transformar los datos para mejorarlos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">No debe: crear información o duplicar, eliminar información relevante</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Mapeo de datos</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Útil cuando se tiene que migrar datos o integrar datos. Pasar de provincia a código, etc..</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Limpieza de datos</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Eliminar los datos nulos o no validos</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Otras tranformaciones:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Cambios de codificación (unificar criterios)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xml)

# The anchor paragraph (originally empty, inheriting the bullet formatting of
# "Kitchen: ejectuar un job") got pushed after the newly inserted content;
# remove it together with the paragraph mark that precedes it.
$beforeAnchor = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$cleanupRange = $d.Range($beforeAnchor.Range.End - 1, $d.Content.End)
$cleanupRange.Delete()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
